# Update the cryptos list: refresh Price (D) and Volume(1h) (E) figures,
# and shift coin rows 40-51 down by one to insert "Frax" at row 40
# (dropping "Chiliz" off the bottom of the list), per the
# "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Many Price values (e.g. "1.005") are digit strings that Excel would
# otherwise auto-convert to numbers. Force the cell to Text first so the
# literal string (incl. trailing zeros) round-trips exactly, then restore
# the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-PlainValue 'D2' '28.261.87'
Set-PlainValue 'E2' '  -6.17%  '
Set-PlainValue 'D3' '1.839.48'
Set-PlainValue 'E3' '  -5.67%  '
Set-TextValue 'D4' '1.005'
Set-PlainValue 'E4' '  -0.45%  '
Set-TextValue 'D5' '332.71'
Set-PlainValue 'E5' '  +1.39%  '
Set-PlainValue 'E6' '  -0.30%  '
Set-TextValue 'D7' '0.4619'
Set-PlainValue 'E7' '  -4.91%  '
Set-PlainValue 'E8' '  -5.65%  '
Set-TextValue 'D9' '46.03'
Set-PlainValue 'E9' '  -3.52%  '
Set-TextValue 'D10' '0.07873'
Set-PlainValue 'E10' '  -4.49%  '
Set-TextValue 'D11' '0.9694'
Set-PlainValue 'E11' '  -5.01%  '
Set-TextValue 'D12' '22.05'
Set-PlainValue 'E12' '  -8.23%  '
Set-PlainValue 'D13' '1.945.62'
Set-PlainValue 'E13' '  +1.84%  '
Set-TextValue 'D14' '5.801'
Set-PlainValue 'E14' '  -5.36%  '
Set-TextValue 'D15' '6.942'
Set-PlainValue 'E15' '  -5.53%  '
Set-TextValue 'D16' '0.06891'
Set-PlainValue 'E16' '  +0.30%  '
Set-TextValue 'D17' '1.006'
Set-PlainValue 'E17' '  -0.45%  '
Set-TextValue 'D18' '87.36'
Set-PlainValue 'E18' '  -5.12%  '
Set-TextValue 'D19' '0.000009976'
Set-PlainValue 'E19' '  -4.18%  '
Set-PlainValue 'E20' '  -5.44%  '
Set-TextValue 'D21' '1.006'
Set-PlainValue 'E21' '  -0.24%  '
Set-PlainValue 'D22' '28.335.19'
Set-PlainValue 'E22' '  -5.86%  '
Set-PlainValue 'E23' '  -5.94%  '
Set-TextValue 'D24' '11.14'
Set-PlainValue 'E24' '  -7.32%  '
Set-TextValue 'D25' '2.171'
Set-PlainValue 'E25' '  -1.40%  '
Set-PlainValue 'D26' '2.098.85'
Set-PlainValue 'E26' '  -2.02%  '
Set-TextValue 'D27' '153.59'
Set-PlainValue 'E27' '  -2.23%  '
Set-TextValue 'D28' '19.33'
Set-PlainValue 'E28' '  -4.26%  '
Set-TextValue 'D29' '5.929'
Set-PlainValue 'E29' '  -10.05%  '
Set-TextValue 'D30' '1.982'
Set-PlainValue 'E30' '  -6.64%  '
Set-TextValue 'D31' '117.18'
Set-PlainValue 'E31' '  -3.57%  '
Set-TextValue 'D32' '0.9508'
Set-PlainValue 'E32' '  -7.39%  '
Set-TextValue 'D33' '0.09350'
Set-PlainValue 'E33' '  -3.10%  '
Set-TextValue 'D34' '5.333'
Set-PlainValue 'E34' '  -5.68%  '
Set-TextValue 'D35' '3.462'
Set-PlainValue 'E35' '  -2.57%  '
Set-TextValue 'D36' '1.330'
Set-PlainValue 'E36' '  -7.11%  '
Set-TextValue 'D37' '0.06058'
Set-PlainValue 'E37' '  -7.44%  '
Set-TextValue 'D38' '0.02180'
Set-PlainValue 'E38' '  -5.63%  '
Set-TextValue 'D39' '1.154'
Set-PlainValue 'E39' '  -6.32%  '
Set-PlainValue 'B40' 'Frax'
Set-PlainValue 'C40' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D40' '1.005'
Set-PlainValue 'E40' '  -0.26%  '
Set-PlainValue 'B41' 'TheSandbox'
Set-PlainValue 'C41' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D41' '0.5657'
Set-PlainValue 'E41' '  -5.62%  '
Set-PlainValue 'B42' 'FraxShare'
Set-PlainValue 'C42' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D42' '7.610'
Set-PlainValue 'E42' '  -4.99%  '
Set-PlainValue 'B43' 'Aptos'
Set-PlainValue 'C43' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D43' '10.02'
Set-PlainValue 'E43' '  -7.08%  '
Set-PlainValue 'B44' 'Algorand'
Set-PlainValue 'C44' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue 'D44' '0.1789'
Set-PlainValue 'E44' '  -3.92%  '
Set-PlainValue 'B45' 'RenderToken'
Set-PlainValue 'C45' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D45' '2.391'
Set-PlainValue 'E45' '  -5.76%  '
Set-PlainValue 'B46' 'WEMIXToken'
Set-PlainValue 'C46' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D46' '1.224'
Set-PlainValue 'E46' '  -1.95%  '
Set-TextValue 'D47' '11.68'
Set-PlainValue 'E47' '  -6.68%  '
Set-PlainValue 'B48' 'Decentraland'
Set-PlainValue 'C48' 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue 'D48' '0.5326'
Set-PlainValue 'E48' '  -4.88%  '
Set-PlainValue 'B49' 'Cronos'
Set-PlainValue 'C49' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D49' '0.07037'
Set-PlainValue 'E49' '  -7.04%  '
Set-PlainValue 'B50' 'NEARProtocol'
Set-PlainValue 'C50' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D50' '1.853'
Set-PlainValue 'E50' '  -7.31%  '
Set-PlainValue 'B51' 'Quant'
Set-PlainValue 'C51' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D51' '112.95'
Set-PlainValue 'E51' '  -4.03%  '
